$d = $word.ActiveDocument

# 1. "Your Name" -> "Nhat Tran" (Author paragraph)
$d.Content.Find.Execute('Your Name', $true, $false, $false, $false, $false, $true, 1, $false, 'Nhat Tran', 2)

# 2. Remove the "System Information" Heading1 paragraph and the
#    "This document shows system details from the EC2 instance." paragraph.
#    (Range.Text includes the trailing paragraph mark, so trim it before
#    comparing.)
$pHeading = $null
$pIntro = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd([char]13, [char]7, [char]11)
    if ($t -eq 'System Information') { $pHeading = $p }
    if ($t -eq 'This document shows system details from the EC2 instance.') { $pIntro = $p }
}
$r = $d.Range($pHeading.Range.Start, $pIntro.Range.End)
$r.Delete()

# 3. The paragraph that used to follow ("import platform ... import psutil")
#    now becomes the first paragraph of the body and picks up the
#    "First Paragraph" style instead of "Body Text".
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd([char]13, [char]7, [char]11)
    if ($t -eq 'import platform import psutil') {
        $p.Style = 'First Paragraph'
        break
    }
}

# 4. Remove everything from the "Press Ctrl+O..." bullet list through the
#    end of the "Step 3" section (including the horizontal rule, the
#    Heading3 paragraph and the closing code block paragraph).
$pFirstBullet = $null
$pLast = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd([char]13, [char]7, [char]11)
    if ($t -eq 'Press Ctrl+O, then Enter to save') { $pFirstBullet = $p }
    if ($t -eq '```bash quarto render system_info.qmd') { $pLast = $p }
}
$r2 = $d.Range($pFirstBullet.Range.Start, $pLast.Range.End)
$r2.Delete()
